$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 6).Value = 2.08
$ws.Cells.Item(2, 7).Value = 2.4
$ws.Cells.Item(2, 8).Value = 2.86
$ws.Cells.Item(2, 9).Value = 3.4
$ws.Cells.Item(2, 10).Value = 3.8
$ws.Cells.Item(2, 11).Value = 5.4
$ws.Cells.Item(2, 12).Value = 1.24
$ws.Cells.Item(2, 14).Value = 6.2
$ws.Cells.Item(2, 15).Value = 1.14
$ws.Cells.Item(2, 16).Value = 2.64
$ws.Cells.Item(2, 18).Value = 1.78
$ws.Cells.Item(2, 19).Value = 2.04
$ws.Cells.Item(2, 22).Value = 1.41
$ws.Cells.Item(2, 23).Value = 1.72
$ws.Cells.Item(2, 40).Value = 85

# Row 3
$ws.Cells.Item(3, 6).Value = 1.09
$ws.Cells.Item(3, 7).Value = 600
$ws.Cells.Item(3, 8).Value = 1.06
$ws.Cells.Item(3, 10).Value = 1.09
$ws.Cells.Item(3, 11).Value = 500
$ws.Cells.Item(3, 12).Value = 1.27
$ws.Cells.Item(3, 14).Value = 1.65
$ws.Cells.Item(3, 18).Value = 1.24
$ws.Cells.Item(3, 19).Value = 1.05
$ws.Cells.Item(3, 20).Value = 1.04
$ws.Cells.Item(3, 21).Value = 1.04

# Row 4
$ws.Cells.Item(4, 14).Value = 1.1
$ws.Cells.Item(4, 20).Value = 1.04
$ws.Cells.Item(4, 21).Value = 1.04

# Row 5
$ws.Cells.Item(5, 7).Value = 1.87
$ws.Cells.Item(5, 8).Value = 4.2
$ws.Cells.Item(5, 23).Value = 2.14
$ws.Cells.Item(5, 35).Value = 170
$ws.Cells.Item(5, 37).Value = 17

# Row 6
$ws.Cells.Item(6, 6).Value = 1.45
$ws.Cells.Item(6, 7).Value = 1.46
$ws.Cells.Item(6, 12).Value = 1.32
$ws.Cells.Item(6, 14).Value = 5.4
$ws.Cells.Item(6, 16).Value = 2.48
$ws.Cells.Item(6, 17).Value = 1.65
$ws.Cells.Item(6, 18).Value = 1.6
$ws.Cells.Item(6, 19).Value = 2.62
$ws.Cells.Item(6, 20).Value = 1.86
$ws.Cells.Item(6, 23).Value = 3.15
$ws.Cells.Item(6, 25).Value = 32
$ws.Cells.Item(6, 27).Value = 260
$ws.Cells.Item(6, 28).Value = 10
$ws.Cells.Item(6, 30).Value = 29

# Row 7
$ws.Cells.Item(7, 6).Value = 2.2
$ws.Cells.Item(7, 7).Value = 2.22
$ws.Cells.Item(7, 11).Value = 3.95
$ws.Cells.Item(7, 19).Value = 2.66
$ws.Cells.Item(7, 20).Value = 1.6
$ws.Cells.Item(7, 23).Value = 1.82
$ws.Cells.Item(7, 24).Value = 20
$ws.Cells.Item(7, 25).Value = 17.5
$ws.Cells.Item(7, 29).Value = 8.6
$ws.Cells.Item(7, 35).Value = 38
$ws.Cells.Item(7, 36).Value = 28
$ws.Cells.Item(7, 39).Value = 65
$ws.Cells.Item(7, 40).Value = 11.5

# Row 8
$ws.Cells.Item(8, 6).Value = 1.51
$ws.Cells.Item(8, 7).Value = 1.61
$ws.Cells.Item(8, 8).Value = 5.3
$ws.Cells.Item(8, 10).Value = 4.9
$ws.Cells.Item(8, 11).Value = 5.9
$ws.Cells.Item(8, 13).Value = 1.01
$ws.Cells.Item(8, 14).Value = 7.4
$ws.Cells.Item(8, 15).Value = 1.12
$ws.Cells.Item(8, 16).Value = 3.2
$ws.Cells.Item(8, 17).Value = 1.35
$ws.Cells.Item(8, 18).Value = 2
$ws.Cells.Item(8, 19).Value = 1.9
$ws.Cells.Item(8, 20).Value = 1.5
$ws.Cells.Item(8, 21).Value = 2.56
$ws.Cells.Item(8, 23).Value = 2.62
$ws.Cells.Item(8, 24).Value = 42
$ws.Cells.Item(8, 25).Value = 40
$ws.Cells.Item(8, 28).Value = 17
$ws.Cells.Item(8, 31).Value = 70
$ws.Cells.Item(8, 39).Value = 260
$ws.Cells.Item(8, 40).Value = 29
$ws.Cells.Item(8, 41).Value = 46

# Row 9
$ws.Cells.Item(9, 6).Value = 2.64
$ws.Cells.Item(9, 7).Value = 3.05
$ws.Cells.Item(9, 8).Value = 2.44
$ws.Cells.Item(9, 9).Value = 2.62
$ws.Cells.Item(9, 10).Value = 3.85
$ws.Cells.Item(9, 12).Value = 1.26
$ws.Cells.Item(9, 14).Value = 5.5
$ws.Cells.Item(9, 15).Value = 1.18
$ws.Cells.Item(9, 16).Value = 2.54
$ws.Cells.Item(9, 17).Value = 1.43
$ws.Cells.Item(9, 18).Value = 1.65
$ws.Cells.Item(9, 19).Value = 2.1
$ws.Cells.Item(9, 20).Value = 1.04
$ws.Cells.Item(9, 21).Value = 2.44
$ws.Cells.Item(9, 22).Value = 1.62
$ws.Cells.Item(9, 24).Value = 1000
$ws.Cells.Item(9, 25).Value = 1000
$ws.Cells.Item(9, 26).Value = 1000
$ws.Cells.Item(9, 27).Value = 1000
$ws.Cells.Item(9, 28).Value = 1000
$ws.Cells.Item(9, 29).Value = 1000
$ws.Cells.Item(9, 30).Value = 1000
$ws.Cells.Item(9, 31).Value = 1000
$ws.Cells.Item(9, 32).Value = 1000
$ws.Cells.Item(9, 33).Value = 1000
$ws.Cells.Item(9, 34).Value = 1000
$ws.Cells.Item(9, 35).Value = 1000
$ws.Cells.Item(9, 37).Value = 1000
$ws.Cells.Item(9, 38).Value = 1000
$ws.Cells.Item(9, 39).Value = 1000
$ws.Cells.Item(9, 40).Value = 1000
$ws.Cells.Item(9, 41).Value = 1000

# Row 10
$ws.Cells.Item(10, 9).Value = 8
$ws.Cells.Item(10, 11).Value = 5.1
$ws.Cells.Item(10, 22).Value = 1.15

# Row 11
$ws.Cells.Item(11, 16).Value = 2.32
$ws.Cells.Item(11, 19).Value = 2.84
$ws.Cells.Item(11, 28).Value = 15.5
$ws.Cells.Item(11, 35).Value = 30
$ws.Cells.Item(11, 41).Value = 14.5

# Row 12
$ws.Cells.Item(12, 12).Value = 1.3
$ws.Cells.Item(12, 14).Value = 5.5
$ws.Cells.Item(12, 16).Value = 2.54
$ws.Cells.Item(12, 17).Value = 1.63
$ws.Cells.Item(12, 18).Value = 1.61
$ws.Cells.Item(12, 19).Value = 2.6
$ws.Cells.Item(12, 25).Value = 10.5
$ws.Cells.Item(12, 28).Value = 34
$ws.Cells.Item(12, 38).Value = 90
$ws.Cells.Item(12, 41).Value = 5.3

# Row 13
$ws.Cells.Item(13, 7).Value = 2.98
$ws.Cells.Item(13, 17).Value = 1.53
$ws.Cells.Item(13, 29).Value = 10
$ws.Cells.Item(13, 31).Value = 21
$ws.Cells.Item(13, 38).Value = 28

# Row 14
$ws.Cells.Item(14, 8).Value = 21
$ws.Cells.Item(14, 16).Value = 2.9
$ws.Cells.Item(14, 18).Value = 1.77
$ws.Cells.Item(14, 20).Value = 2.36
$ws.Cells.Item(14, 21).Value = 1.7
$ws.Cells.Item(14, 31).Value = 440
$ws.Cells.Item(14, 34).Value = 44
$ws.Cells.Item(14, 38).Value = 48

# Row 15
$ws.Cells.Item(15, 8).Value = 1.73
$ws.Cells.Item(15, 9).Value = 1.74
$ws.Cells.Item(15, 14).Value = 4.5
$ws.Cells.Item(15, 16).Value = 2.18
$ws.Cells.Item(15, 20).Value = 1.82
$ws.Cells.Item(15, 22).Value = 2.34
$ws.Cells.Item(15, 27).Value = 16.5
$ws.Cells.Item(15, 30).Value = 9.800000000000001
$ws.Cells.Item(15, 31).Value = 16.5
$ws.Cells.Item(15, 34).Value = 19
$ws.Cells.Item(15, 38).Value = 65
$ws.Cells.Item(15, 39).Value = 95
$ws.Cells.Item(15, 40).Value = 70
$ws.Cells.Item(15, 41).Value = 8.800000000000001

# Row 16
$ws.Cells.Item(16, 8).Value = 12.5
$ws.Cells.Item(16, 9).Value = 13
$ws.Cells.Item(16, 14).Value = 9
$ws.Cells.Item(16, 16).Value = 3.65
$ws.Cells.Item(16, 28).Value = 16
$ws.Cells.Item(16, 32).Value = 11
$ws.Cells.Item(16, 39).Value = 85
$ws.Cells.Item(16, 41).Value = 1000

# Row 17
$ws.Cells.Item(17, 6).Value = 2.5
$ws.Cells.Item(17, 7).Value = 2.52
$ws.Cells.Item(17, 8).Value = 3.25
$ws.Cells.Item(17, 12).Value = 1.49
$ws.Cells.Item(17, 14).Value = 3.4
$ws.Cells.Item(17, 15).Value = 1.4
$ws.Cells.Item(17, 16).Value = 1.79
$ws.Cells.Item(17, 17).Value = 2.22
$ws.Cells.Item(17, 18).Value = 1.3
$ws.Cells.Item(17, 19).Value = 4.2
$ws.Cells.Item(17, 20).Value = 1.89
$ws.Cells.Item(17, 21).Value = 2.04
$ws.Cells.Item(17, 22).Value = 1.43
$ws.Cells.Item(17, 23).Value = 1.65
$ws.Cells.Item(17, 25).Value = 11
$ws.Cells.Item(17, 28).Value = 9.4
$ws.Cells.Item(17, 38).Value = 46
$ws.Cells.Item(17, 39).Value = 120
$ws.Cells.Item(17, 40).Value = 25

# Row 18
$ws.Cells.Item(18, 7).Value = 2.7
$ws.Cells.Item(18, 8).Value = 2.54
$ws.Cells.Item(18, 9).Value = 2.96
$ws.Cells.Item(18, 10).Value = 3.65
$ws.Cells.Item(18, 14).Value = 5.5
$ws.Cells.Item(18, 16).Value = 2.88
$ws.Cells.Item(18, 17).Value = 1.48
$ws.Cells.Item(18, 18).Value = 1.73
$ws.Cells.Item(18, 19).Value = 2.24
$ws.Cells.Item(18, 20).Value = 1.48
$ws.Cells.Item(18, 21).Value = 1.04
$ws.Cells.Item(18, 22).Value = 1.51
$ws.Cells.Item(18, 23).Value = 1.59
$ws.Cells.Item(18, 24).Value = 28
$ws.Cells.Item(18, 25).Value = 19.5
$ws.Cells.Item(18, 26).Value = 25
$ws.Cells.Item(18, 27).Value = 110
$ws.Cells.Item(18, 28).Value = 18.5
$ws.Cells.Item(18, 29).Value = 11.5
$ws.Cells.Item(18, 30).Value = 14.5
$ws.Cells.Item(18, 31).Value = 27
$ws.Cells.Item(18, 32).Value = 23
$ws.Cells.Item(18, 33).Value = 13
$ws.Cells.Item(18, 35).Value = 70
$ws.Cells.Item(18, 37).Value = 24
$ws.Cells.Item(18, 39).Value = 60
$ws.Cells.Item(18, 40).Value = 14
$ws.Cells.Item(18, 41).Value = 16

# Row 19
$ws.Cells.Item(19, 9).Value = 3.2
$ws.Cells.Item(19, 12).Value = 1.42
$ws.Cells.Item(19, 16).Value = 1.82
$ws.Cells.Item(19, 19).Value = 3.6
$ws.Cells.Item(19, 21).Value = 2.08
$ws.Cells.Item(19, 25).Value = 14
$ws.Cells.Item(19, 31).Value = 38
$ws.Cells.Item(19, 35).Value = 50

# Row 20
$ws.Cells.Item(20, 6).Value = 1.55
$ws.Cells.Item(20, 9).Value = 6.8
$ws.Cells.Item(20, 10).Value = 1.18
$ws.Cells.Item(20, 19).Value = 1.31
$ws.Cells.Item(20, 20).Value = 1.04
$ws.Cells.Item(20, 21).Value = 1.04
$ws.Cells.Item(20, 22).Value = 1.17
